# Add 2022-Q4 data
# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" positioned right after "总计"
#    and before "2022-Q3" (i.e. directly before the first quarter sheet).
# 2) Populate it with the new fund-holding rows for that quarter.
# 3) Insert a new row into the "总计" summary sheet for 2022-Q4 and
#    renumber the index column of the rows that shifted down.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" sheet in the right position ------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"

# Re-fetch by name (handle returned by Add() becomes unreliable once
# the sheet collection is reordered) and move it before "2022-Q3".
$q4 = $wb.Worksheets.Item("2022-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4.Move($q3)

$ws = $wb.Worksheets.Item("2022-Q4")
$src = $wb.Worksheets.Item("2022-Q3")

# Carry over the existing header / index-column formatting so the new
# sheet visually matches its siblings.
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2:A4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# --- Header row ------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- Data rows (text-like numeric fields keep a leading apostrophe so
#     the stored type stays text, matching the source data) ----------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'010149"
$ws.Range("C2").Value = "浙商智选经济动能混合C"
$ws.Range("D2").Value = "'4.34"
$ws.Range("E2").Value = "'87.31"
$ws.Range("F2").Value = "'5.56"
$ws.Range("G2").Value = "'0.2413"
$ws.Range("H2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'001449"
$ws.Range("C3").Value = "华商双驱优选灵活配置混合"
$ws.Range("D3").Value = "'2.26"
$ws.Range("E3").Value = "'77.71"
$ws.Range("F3").Value = "'3.65"
$ws.Range("G3").Value = "'0.0825"
$ws.Range("H3").Value = 6

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'010148"
$ws.Range("C4").Value = "浙商智选经济动能混合A"
$ws.Range("D4").Value = "'0.62"
$ws.Range("E4").Value = "'87.31"
$ws.Range("F4").Value = "'5.56"
$ws.Range("G4").Value = "'0.0345"
$ws.Range("H4").Value = 3

# --- 2. Update the "总计" summary sheet ------------------------------
$total = $wb.Worksheets.Item("总计")

# Insert a blank row at row 2 for the new quarter; this pushes the
# existing 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q2 rows down by one.
$total.Rows("2:2").Insert()

# Pick up plain (unstyled) formatting from the row directly below for
# the new data row, and the bold/bordered index-column style for A2.
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.36

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
